# Update scripts with new TPM-derived values for the Plg-F2r LR-pair sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs -> ECs)
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 13.83830466666667
$ws.Range("N2").Value = 41.514914
$ws.Range("O2").Value = 0.1568893148900199
$ws.Range("P2").Value = 0.1568893148900199
$ws.Range("Q2").Value = 0.1955260194035556
$ws.Range("R2").Value = 1.759734174632
$ws.Range("S2").Value = 0.1568893148900199
$ws.Range("T2").Value = 0.1568893148900199

# Row 3 (FAPs -> FAPs)
$ws.Range("O3").Value = 0.5509859018285573
$ws.Range("P3").Value = 0.5509859018285573
$ws.Range("Q3").Value = 0.6866757000471112
$ws.Range("R3").Value = 6.180081300424001
$ws.Range("S3").Value = 0.5509859018285573
$ws.Range("T3").Value = 0.5509859018285573

# Row 4 (FAPs -> MuSCs)
$ws.Range("M4").Value = 25.766648
$ws.Range("N4").Value = 77.299944
$ws.Range("O4").Value = 0.2921247832814228
$ws.Range("P4").Value = 0.2921247832814228
$ws.Range("Q4").Value = 0.3640655584746667
$ws.Range("R4").Value = 3.276590026272
$ws.Range("S4").Value = 0.2921247832814228
$ws.Range("T4").Value = 0.2921247832814228
